$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New 20-minute trade row (row 4) appended below the existing trade rows.
# Values mirror the columns: Principle, Start Principle, BuyPrice, SellPrice,
# IsShortSell, Price Change %, Date, Profitable
$ws.Range("A4").Value = 10015.91
$ws.Range("B4").Value = 10039
$ws.Range("C4").Value = 286.39
$ws.Range("D4").Value = 287.04000000000002
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = 0.23

# Copy the Date column's number format from the row above so the new Date
# cell keeps the same style index (m/d/yyyy h:mm) instead of minting a new one.
$ws.Range("G3").Copy()
$ws.Range("G4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G4").Value = 42608.639641203707

$ws.Range("H4").Value = $false

# Column A now needs to fit the wider "10015.91" value (natural best-fit
# width of 9 characters).
$ws.Columns.Item(1).ColumnWidth = 8.1
